{"js": "// The intro paragraph used to read:\n//   ...hereinafter referred to as \"Parties,\" have entered...\n// The comma that sat just inside the closing quotation mark moves to\n// just outside it:\n//   ...hereinafter referred to as \"Parties\", have entered...\nconst body = context.document.body;\n\nconst oldPhrase = 'hereinafter referred to as \"Parties,\" have entered';\nconst newPhrase = 'hereinafter referred to as \"Parties\", have entered';\n\nconst results = body.search(oldPhrase, { matchCase: true });\nresults.load('items');\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(newPhrase, 'Replace');\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The intro paragraph used to read:\n#   ...hereinafter referred to as \"Parties,\" have entered...\n# The comma that sat just inside the closing quotation mark moves to\n# just outside it:\n#   ...hereinafter referred to as \"Parties\", have entered...\n$quote = [char]34\n$oldPhrase = $quote + 'Parties,' + $quote\n$newPhrase = $quote + 'Parties' + $quote + ','\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldPhrase\n$find.Execute()\n\nif ($find.Found) {\n    $rng = $find.Parent\n    $rng.Text = $newPhrase\n}\n"}
